$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Update the stat values in the first column (rows are 1-indexed in Word COM).
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "7161"
$t.Cell(7, 1).Range.Text  = "0.06545"
$t.Cell(8, 1).Range.Text  = "0.07467"
$t.Cell(12, 1).Range.Text = "440.45393"

# The trailing three rows previously held full tab-separated readouts; they
# are collapsed down to the single summary value that now lives elsewhere.
$t.Cell(44, 1).Range.Text = "59.49"
$t.Cell(45, 1).Range.Text = "440.45"
$t.Cell(46, 1).Range.Text = "1087"
